$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param(
        [string]$rowRange,
        [object[]]$values
    )
    $arr = New-Object 'object[,]' 1,$values.Length
    for ($i = 0; $i -lt $values.Length; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $ws.Range($rowRange).Value = $arr
}

# Row 2: ECs -> Vegfc -> Vipr2 -> ECs
Set-RowValues "A2:T2" @(
    "ECs","Vegfc","Vipr2","ECs",
    2,0.6666666666666666,1.164278333333333,3.492835,0.2070506538112546,0.2070506538112546,
    1,0.3333333333333333,0.02669166666666667,0.08007499999999999,0.002946855435681398,0.002946855435681398,
    0.03107652918055556,0.279688762625,0.000610148344645083,0.000610148344645083
)

# Row 3: ECs -> Vegfc -> Vipr2 -> FAPs
Set-RowValues "A3:T3" @(
    "ECs","Vegfc","Vipr2","FAPs",
    2,0.6666666666666666,1.164278333333333,3.492835,0.2070506538112546,0.2070506538112546,
    3,1,4.865925666666667,14.597777,0.5372155916492649,0.5372155916492649,
    5.665291825310556,50.987626427795,0.1112308394885803,0.1112308394885803
)

# Row 4: ECs -> Vegfc -> Vipr2 -> sCs
Set-RowValues "A4:T4" @(
    "ECs","Vegfc","Vipr2","sCs",
    2,0.6666666666666666,1.164278333333333,3.492835,0.2070506538112546,0.2070506538112546,
    3,1,4.165060333333333,12.495181,0.4598375529150536,0.4598375529150536,
    4.84928950312611,43.643605528135,0.09520966597802925,0.09520966597802925
)

# Row 5: FAPs -> Vegfc -> Vipr2 -> ECs
Set-RowValues "A5:T5" @(
    "FAPs","Vegfc","Vipr2","ECs",
    3,1,3.886460333333333,11.659381,0.6911527338349851,0.6911527338349851,
    1,0.3333333333333333,0.02669166666666667,0.08007499999999999,0.002946855435681398,0.002946855435681398,
    0.1037361037305556,0.9336249335749999,0.002036727190587684,0.002036727190587684
)

# Row 6: FAPs -> Vegfc -> Vipr2 -> FAPs
Set-RowValues "A6:T6" @(
    "FAPs","Vegfc","Vipr2","FAPs",
    3,1,3.886460333333333,11.659381,0.6911527338349851,0.6911527338349851,
    3,1,4.865925666666667,14.597777,0.5372155916492649,0.5372155916492649,
    18.91122708844856,170.201043796037,0.3712980248271684,0.3712980248271684
)

# Row 7: FAPs -> Vegfc -> Vipr2 -> sCs
Set-RowValues "A7:T7" @(
    "FAPs","Vegfc","Vipr2","sCs",
    3,1,3.886460333333333,11.659381,0.6911527338349851,0.6911527338349851,
    3,1,4.165060333333333,12.495181,0.4598375529150536,0.4598375529150536,
    16.18734177144011,145.686075942961,0.3178179818172289,0.3178179818172289
)

# Row 8: sCs -> Vegfc -> Vipr2 -> ECs
Set-RowValues "A8:T8" @(
    "sCs","Vegfc","Vipr2","ECs",
    3,1,0.5724183333333334,1.717255,0.1017966123537602,0.1017966123537602,
    1,0.3333333333333333,0.02669166666666667,0.08007499999999999,0.002946855435681398,0.002946855435681398,
    0.01527879934722222,0.137509194125,0.0002999799004486304,0.0002999799004486304
)

# Row 9: sCs -> Vegfc -> Vipr2 -> FAPs
Set-RowValues "A9:T9" @(
    "sCs","Vegfc","Vipr2","FAPs",
    3,1,0.5724183333333334,1.717255,0.1017966123537602,0.1017966123537602,
    3,1,4.865925666666667,14.597777,0.5372155916492649,0.5372155916492649,
    2.785345060237222,25.068105542135,0.05468672733351618,0.05468672733351617
)

# Row 10: sCs -> Vegfc -> Vipr2 -> sCs
Set-RowValues "A10:T10" @(
    "sCs","Vegfc","Vipr2","sCs",
    3,1,0.5724183333333334,1.717255,0.1017966123537602,0.1017966123537602,
    3,1,4.165060333333333,12.495181,0.4598375529150536,0.4598375529150536,
    2.384156894239444,21.457412048155,0.04680990511979542,0.04680990511979542
)
